# Add a "president" column (F) to the approval data sheet and fill it
# with the constant value "Bush, sr" for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new column.
$ws.Range("F1").Value = "president"

# One "Bush, sr" value per existing data row (rows 2 through 159).
$ws.Range("F2:F159").Value = "Bush, sr"

# Match the formatting Excel applies to the other header cells
# (bold font + border via the shared header style) by copying it
# from the neighboring header cell onto the new one, without
# touching its value.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
